$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 175542
$ws.Range("C4").Value = 165515
$ws.Range("C7").Value = 5.71
$ws.Range("C8").Value = 64.53
